$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "REX_DEF" in F1, matching the style of the existing header cells (E1)
$ws.Range("F1").Value = "REX_DEF"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
